$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old contents (4 rows x 5 cols) entirely before rewriting.
$ws.Range("A1:E4").Clear()

# Header row
$ws.Range("A1").Value = "매뉴얼ID"
$ws.Range("B1").Value = "한글제목"
$ws.Range("C1").Value = "한글요약"
$ws.Range("D1").Value = "한글본문"

# First data row
$ws.Range("A2").Value = "example-entry-1"
$ws.Range("B2").Formula = "="""""
$ws.Range("C2").Formula = "="""""
$ws.Range("D2").Formula = "="""""
